$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "NewBvDIDs",
    "FR542051180",
    "DE5110206610",
    "DE5330000056",
    "US142546192L",
    "GB12133117",
    "DE6270139382",
    "DE5110426441",
    "DE5190434206",
    "DE4010000539",
    "DE5030088293",
    "FR180020026",
    "DE8330415355",
    "PL811655650",
    "ESB98338411",
    "DE2151574829",
    "DE5050471540",
    "DE6070009834",
    "GB*908439227",
    "DE5050567293",
    "NO988470937",
    "DE5170178473",
    "DE5030319504",
    "SE5564620887",
    "BM45532R",
    "NO863769132",
    "DE2190001854",
    "DE3170295464",
    "GB01142830",
    "US870698303",
    "LULB201907",
    "DEFEB41239",
    "DE5030456699",
    "DEFEB13116",
    "DE2290291348",
    "FR531330389",
    "DE5110426307",
    "DE2370010537",
    "SE5560362138",
    "ESA08001851",
    "DK10103940",
    "KY*110014084853",
    "CZ29259428",
    "ZA196701303806",
    "ESB82056748",
    "SE5567818314",
    "DE6070491599",
    "DE2210163123",
    "AT9110778886",
    "ESA28119220",
    "GB00102498",
    "LULB195709",
    "DE7170000175",
    "US133413761L",
    "DE6150261280",
    "DE2011364466",
    "IT04732570967",
    "CHCHE112112310",
    "DE2050372206",
    "DE8170000724",
    "DE2150918982",
    "DE3070183858",
    "DE3070407436",
    "DE5050056484",
    "GBGG46130",
    "ESA28157360",
    "DE6150007515",
    "DEFEB15469",
    "PL012149187",
    "CYC307511",
    "IT06312510966",
    "DE7170311222",
    "DE2350060811",
    "IN0008687067",
    "DE8190350014",
    "DE7230282328",
    "DE8170504074",
    "NL53748662",
    "DE5110304100",
    "US150238465L",
    "DE8170725443",
    "KY31377PC",
    "US411724239",
    "EE75019980",
    "DE2151253621",
    "CHCHE109101085",
    "US541163725",
    "DE8070234353",
    "DE2150965540",
    "FR443837331",
    "DE5110066894",
    "DE5070255751",
    "GBLP011863",
    "IT10236451000",
    "DE2250192021",
    "GB04043908",
    "DE2350113248",
    "DE8310297852",
    "CHCHE114669799",
    "AT9150100345",
    "DE2011425555",
    "NL33292246",
    "LULB104637",
    "IT00811720580",
    "DE7110175411",
    "GBGG50009",
    "DE2230133950",
    "KY*110179911376",
    "DE7330620677",
    "NL01178978",
    "DE*929123061",
    "NL34108286",
    "DE2190331916",
    "CZ00260428",
    "DE6290207493",
    "IS6306051210",
    "RU00040778",
    "NO987059699",
    "GB02989338",
    "DE2250063387",
    "DE7370113108",
    "DE2050488817",
    "DE6070524807",
    "US363145972",
    "RU17230282",
    "DEFEB23470",
    "DE4070408675",
    "DE4170123108",
    "DE2250132994",
    "FR552081317",
    "DEFEB13584",
    "US149185809L",
    "DE8290178515",
    "GB03039100",
    "GBOC303117",
    "FR379677636",
    "FR*J00S1080",
    "DE7330530056",
    "IE008610",
    "DEFEI1007130",
    "DE5110216866",
    "FR542008677",
    "DE8170600335",
    "DE*820975510",
    "DE5110001494",
    "CH0000049173",
    "HU10625790",
    "DE7010228841",
    "DE8190010796",
    "DE2010198197",
    "CZ41692918",
    "DEFEB13109",
    "NO912230252",
    "DE2030006515",
    "FR380415125",
    "DE2010378297",
    "ESG08171118",
    "DEFEB13306",
    "FI06353667",
    "US341843785",
    "GB03586615",
    "DE*J00S1100",
    "GB02366963",
    "DEFEB13172",
    "DE7330029432",
    "GB02020410",
    "GB01833679",
    "DEFEB48901",
    "DE6070002888",
    "DEFEI1007380",
    "GB03690065",
    "FI02459755"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
